# Spotify updated, Code caught up!
# Update header "NEW PLAYS" date, refresh the first track's stats, and
# remove the four stale tracks (rows 3-6) that are no longer part of the
# current top list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: new-plays-as-of date
$ws.Range("B1").Value = "NEW PLAYS September 01, 2025"

# Row 2 (Sundarakalebara) refreshed play count.
# Force text formatting first so the comma-formatted number stays a
# literal string (matching the source data) instead of being parsed as
# a numeric value, then drop back to the default "Normal" style so no
# stray number-format style lingers on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "119,494"
$ws.Range("B2").Style = "Normal"

# Release date reformatted from ISO to long form.
$ws.Range("C2").Value = "August 9, 2024"

# Artist credit trimmed to the primary artist.
$ws.Range("D2").Value = "Bineetha Ranjith"

# Drop the other four tracks entirely - they fell out of the list.
$ws.Rows("3:6").Delete()
